$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Commands" sheet: insert the new Artisan Command keepON(<bool>) right
# after the existing keyboard(<bool>) row (old row 98). This pushes the
# RC Command / YOCTOPUCE / WebSocket Command rows down by one row and
# effectively creates one brand-new trailing row, exactly like the
# original commit's change to doc/help_dialogs/Input_files/eventsliders.xlsx.
# ---------------------------------------------------------------------
$wsCommands = $wb.Worksheets.Item("Commands")

$wsCommands.Rows.Item(99).Insert() | Out-Null
$wsCommands.Range("B99").Value2 = "keepON(<bool>)"
$wsCommands.Range("C99").Value2 = "enables/disables the Keep ON flag"
$wsCommands.Rows.Item(99).RowHeight = 13.8

# Update the view/selection to reflect the newly inserted row (B99:C99
# replaces the old B20:C20 selection).
$wsCommands.Range("B99:C99").Select() | Out-Null

# ---------------------------------------------------------------------
# "Sliders" sheet: only the selection bookkeeping changes (it mirrors
# the updated row reference on the Commands sheet), no data changes.
# ---------------------------------------------------------------------
$wsSliders = $wb.Worksheets.Item("Sliders")
$wsSliders.Range("B99:C99").Select() | Out-Null
$wsSliders.Range("B6").Activate() | Out-Null

$wsCommands.Select() | Out-Null
